$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text / percentage / name / URL cells - safe to set directly
$ws.Range('E2').Value = '  +2.13%  '
$ws.Range('E3').Value = '  -0.05%  '
$ws.Range('E4').Value = '  -0.59%  '
$ws.Range('E5').Value = '  -0.27%  '
$ws.Range('E6').Value = '  -0.59%  '
$ws.Range('E7').Value = '  +1.58%  '
$ws.Range('E8').Value = '  +0.99%  '
$ws.Range('E9').Value = '  +3.63%  '
$ws.Range('E10').Value = '  +0.88%  '
$ws.Range('E11').Value = '  +0.52%  '
$ws.Range('E12').Value = '  -0.61%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('E13').Value = '  +2.02%  '
$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('E14').Value = '  +3.84%  '
$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('E15').Value = '  +2.52%  '
$ws.Range('E16').Value = '  +1.87%  '
$ws.Range('E17').Value = '  +0.28%  '
$ws.Range('E18').Value = '  -0.61%  '
$ws.Range('E19').Value = '  +1.84%  '
$ws.Range('E20').Value = '  +0.64%  '
$ws.Range('E21').Value = '  +2.27%  '
$ws.Range('E22').Value = '  -0.61%  '
$ws.Range('E23').Value = '  +2.21%  '
$ws.Range('E24').Value = '  +0.54%  '
$ws.Range('E25').Value = '  -0.27%  '
$ws.Range('E26').Value = '  -0.34%  '
$ws.Range('E27').Value = '  +1.19%  '
$ws.Range('E28').Value = '  +7.85%  '
$ws.Range('E29').Value = '  +0.05%  '
$ws.Range('E30').Value = '  +0.35%  '
$ws.Range('E31').Value = '  -0.40%  '
$ws.Range('E32').Value = '  +0.12%  '
$ws.Range('E33').Value = '  -0.26%  '
$ws.Range('E34').Value = '  +1.66%  '
$ws.Range('E35').Value = '  -1.56%  '
$ws.Range('E36').Value = '  -0.55%  '
$ws.Range('E37').Value = '  +2.12%  '
$ws.Range('E38').Value = '  +1.38%  '
$ws.Range('E39').Value = '  +1.77%  '
$ws.Range('E40').Value = '  +0.97%  '
$ws.Range('E41').Value = '  +1.83%  '
$ws.Range('E42').Value = '  -0.24%  '
$ws.Range('E43').Value = '  +1.43%  '
$ws.Range('E44').Value = '  +0.25%  '
$ws.Range('E45').Value = '  +0.94%  '
$ws.Range('E46').Value = '  -0.55%  '
$ws.Range('E47').Value = '  +19.34%  '
$ws.Range('E48').Value = '  +0.28%  '
$ws.Range('E49').Value = '  -1.10%  '
$ws.Range('E50').Value = '  +9.78%  '
$ws.Range('E51').Value = '  -0.20%  '

# Price column (D) - force text to preserve exact formatting (avoid Excel numeric auto-conversion)
$priceCells = $ws.Range("D2:D51")
$priceCells.NumberFormat = "@"
$ws.Range('D2').Value = '30.371.99'
$ws.Range('D3').Value = '2.094.68'
$ws.Range('D4').Value = '1.004'
$ws.Range('D5').Value = '343.08'
$ws.Range('D8').Value = '0.4424'
$ws.Range('D9').Value = '54.58'
$ws.Range('D10').Value = '0.09346'
$ws.Range('D11').Value = '1.169'
$ws.Range('D12').Value = '24.74'
$ws.Range('D13').Value = '2.145.89'
$ws.Range('D14').Value = '8.576'
$ws.Range('D15').Value = '6.922'
$ws.Range('D16').Value = '101.40'
$ws.Range('D20').Value = '0.06688'
$ws.Range('D21').Value = '6.332'
$ws.Range('D23').Value = '30.421.31'
$ws.Range('D24').Value = '12.55'
$ws.Range('D25').Value = '2.308'
$ws.Range('D26').Value = '21.85'
$ws.Range('D27').Value = '163.20'
$ws.Range('D28').Value = '6.777'
$ws.Range('D29').Value = '2.511'
$ws.Range('D30').Value = '133.45'
$ws.Range('D31').Value = '1.137'
$ws.Range('D32').Value = '0.1049'
$ws.Range('D33').Value = '1.645'
$ws.Range('D34').Value = '6.271'
$ws.Range('D35').Value = '3.874'
$ws.Range('D36').Value = '10.16'
$ws.Range('D37').Value = '0.02634'
$ws.Range('D38').Value = '0.06812'
$ws.Range('D39').Value = '0.7008'
$ws.Range('D40').Value = '12.58'
$ws.Range('D41').Value = '1.343'
$ws.Range('D42').Value = '0.2218'
$ws.Range('D43').Value = '0.6866'
$ws.Range('D44').Value = '14.37'
$ws.Range('D45').Value = '2.342'
$ws.Range('D46').Value = '1.002'
$ws.Range('D47').Value = '1.385'
$ws.Range('D49').Value = '0.00000000354'
$ws.Range('D50').Value = '1.232'
$ws.Range('D51').Value = '1.218'
$priceCells.Style = "Normal"
